$d = $word.ActiveDocument

# Update the date title
$d.Content.Find.Execute("2024-10-05 Saturday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-10-06 Sunday", 2)

$t = $d.Tables(1)

# Row 1 (problem set 1)
$t.Cell(1,1).Range.Text = "32÷7="
$t.Cell(1,2).Range.Text = "66÷4="
$t.Cell(1,3).Range.Text = "44÷4="
$t.Cell(1,4).Range.Text = "18÷2="
$t.Cell(1,5).Range.Text = "36÷7="

# Row 5 (problem set 2)
$t.Cell(5,1).Range.Text = "86÷3="
$t.Cell(5,2).Range.Text = "79÷9="
$t.Cell(5,3).Range.Text = "73÷9="
$t.Cell(5,4).Range.Text = "51÷8="
$t.Cell(5,5).Range.Text = "40÷9="

# Row 9 (problem set 3)
$t.Cell(9,1).Range.Text = "56÷4="
$t.Cell(9,2).Range.Text = "19÷9="
$t.Cell(9,3).Range.Text = "48÷7="
$t.Cell(9,4).Range.Text = "86÷7="
$t.Cell(9,5).Range.Text = "27÷7="

# Row 13 (problem set 4)
$t.Cell(13,1).Range.Text = "59÷3="
$t.Cell(13,2).Range.Text = "84÷4="
$t.Cell(13,3).Range.Text = "99÷7="
$t.Cell(13,4).Range.Text = "83÷4="
$t.Cell(13,5).Range.Text = "67÷3="

# Row 17 (problem set 5) - two leading problems dropped, values shift left,
# and two new problems appended at the end
$t.Cell(17,1).Range.Text = "10÷7="
$t.Cell(17,2).Range.Text = "18÷4="
$t.Cell(17,3).Range.Text = "57÷8="
$t.Cell(17,4).Range.Text = "94÷3="
$t.Cell(17,5).Range.Text = "86÷5="
